$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "clave" header to "contraseña"
$ws.Range("B1").Value = "contraseña"

# Update the selected cell to B2
$ws.Range("B2").Select()
